$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.782.29"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.895.85"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'0.7608"
$ws.Range("E5").Value = "  +3.76%  "
$ws.Range("D6").Value = "'239.99"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.3034"
$ws.Range("E8").Value = "  -1.83%  "
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("D10").Value = "'0.06803"
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("D12").Value = "1.890.44"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "'0.7345"
$ws.Range("E13").Value = "  -4.66%  "
$ws.Range("D14").Value = "'5.147"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "'90.73"
$ws.Range("E15").Value = "  -0.76%  "
$ws.Range("D16").Value = "29.780.64"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "'13.78"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "'5.892"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").Value = "'240.92"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "'0.000007685"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "2.134.53"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'6.887"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").Value = "'166.50"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").Value = "'9.203"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'18.60"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "'0.1283"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").Value = "'2.019"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +3.15%  "
$ws.Range("D31").Value = "'1.513"
$ws.Range("D32").Value = "'4.251"
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").Value = "'4.058"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "'0.05188"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").Value = "'1.244"
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("D36").Value = "'0.7224"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").Value = "'2.714"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").Value = "'0.01912"
$ws.Range("D39").Value = "'2.772"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").Value = "'6.135"
$ws.Range("E40").Value = "  -2.59%  "
$ws.Range("D41").Value = "'0.4386"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "'71.83"
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'1.876"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("D45").Value = "'0.8269"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").Value = "'7.593"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "'99.48"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'9.689"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").Value = "2.039.37"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'36.00"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").Value = "'0.05919"
$ws.Range("E51").Value = "  -0.44%  "